# CSE316 scoresheet correction
# - swaps in the correct team's names (Sharanya Kataru / Fahim Jawad)
# - regrades two items that were mis-scored (Use-Case 5 dup-account test,
#   Use-Case 3 network-error test) and adjusts two admin-profile items to
#   half credit with notes
# - replaces the "no zero users" comment with the correct note
# - removes the now-unneeded "Misc. Penalties" line (row 170) and fixes
#   the subtotal formula that referenced it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team / student info (header block) ---
$ws.Range("B4").Value = "sharanya and fahim Finale"
$ws.Range("B5").Value = "Fahim Jawad"
$ws.Range("B6").Value = "Sharanya Kataru"

# --- Use-Case 05 (row 27): duplicate-account test now passes ---
$ws.Range("E27").Value = 1
$ws.Range("G27").Value = ""

# --- Use-Case 03 (row 38): network-error logout test now passes ---
$ws.Range("E38").Value = 1
$ws.Range("G38").Value = ""

# --- Use-Case 19 admin listing (row 141): half credit + note ---
$ws.Range("E141").Value = 0.5
$ws.Range("G141").Value = "not default listing"

# --- Use-Case 19 delete user (row 144): half credit + note ---
$ws.Range("E144").Value = 0.5
$ws.Range("G144").Value = "crash on delete of newly created user, failed to remove communities user had created"

# --- Use-Case 19 no-users message (row 146): corrected comment ---
$ws.Range("G146").Value = "can't empty listing since admin included, no message about 0 regular users"

# --- Column G no longer needs to be as wide now that the long comment is gone ---
$ws.Columns.Item(7).ColumnWidth = 28.65

# --- Drop the old "Misc. Penalties" line (row 170) entirely ---
$ws.Rows.Item(170).Delete()

# Fix the grand-total formula: it used to add the (now deleted) penalty row
$ws.Range("B168").Formula = "=SUM(B164,B166)"
